# Auto-generated Excel COM-interop script to update cryptos list values
# (Sun Apr 16 09:58:51 UTC 2023 GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    # Force text number-format first so numeric-looking strings
    # (e.g. '0.5268') are stored as text, matching the source data,
    # then reset the style back to "Normal" so no stray per-cell
    # style index is left behind (keeps style 0 on these cells).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextCell $ws "D2" "30.652.89"
Set-TextCell $ws "E2" "  +0.49%  "
Set-TextCell $ws "D3" "2.115.56"
Set-TextCell $ws "E3" "  +0.22%  "
Set-TextCell $ws "E4" "  +1.35%  "
Set-TextCell $ws "D5" "337.99"
Set-TextCell $ws "E5" "  +1.15%  "
Set-TextCell $ws "E6" "  +1.26%  "
Set-TextCell $ws "D7" "0.5268"
Set-TextCell $ws "E7" "  +0.75%  "
Set-TextCell $ws "D8" "0.4542"
Set-TextCell $ws "E8" "  +0.95%  "
Set-TextCell $ws "D9" "53.81"
Set-TextCell $ws "E9" "  +0.72%  "
Set-TextCell $ws "D10" "0.09087"
Set-TextCell $ws "E10" "  +1.12%  "
Set-TextCell $ws "D11" "1.172"
Set-TextCell $ws "E11" "  +0.68%  "
Set-TextCell $ws "D12" "24.42"
Set-TextCell $ws "E12" "  -0.24%  "
Set-TextCell $ws "D13" "2.124.71"
Set-TextCell $ws "E13" "  +1.09%  "
Set-TextCell $ws "D14" "6.827"
Set-TextCell $ws "E14" "  +0.56%  "
Set-TextCell $ws "D15" "8.079"
Set-TextCell $ws "E15" "  +3.15%  "
Set-TextCell $ws "D16" "98.39"
Set-TextCell $ws "E16" "  +1.86%  "
Set-TextCell $ws "D17" "0.00001167"
Set-TextCell $ws "E17" "  +3.66%  "
Set-TextCell $ws "D18" "1.015"
Set-TextCell $ws "E18" "  +1.24%  "
Set-TextCell $ws "D19" "0.06698"
Set-TextCell $ws "E19" "  +1.32%  "
Set-TextCell $ws "D20" "19.51"
Set-TextCell $ws "E20" "  +0.95%  "
Set-TextCell $ws "D22" "6.445"
Set-TextCell $ws "E22" "  +2.19%  "
Set-TextCell $ws "D23" "30.729.79"
Set-TextCell $ws "E23" "  +0.58%  "
Set-TextCell $ws "D24" "12.92"
Set-TextCell $ws "E24" "  +4.51%  "
Set-TextCell $ws "D25" "2.378"
Set-TextCell $ws "E25" "  +1.58%  "
Set-TextCell $ws "D26" "2.368.37"
Set-TextCell $ws "E26" "  +0.81%  "
Set-TextCell $ws "D27" "22.42"
Set-TextCell $ws "E27" "  +0.09%  "
Set-TextCell $ws "D28" "165.49"
Set-TextCell $ws "E28" "  +0.98%  "
Set-TextCell $ws "D29" "2.541"
Set-TextCell $ws "E29" "  -1.74%  "
Set-TextCell $ws "D30" "134.95"
Set-TextCell $ws "E30" "  +1.53%  "
Set-TextCell $ws "D31" "1.203"
Set-TextCell $ws "E31" "  -0.06%  "
Set-TextCell $ws "D32" "0.1077"
Set-TextCell $ws "E32" "  +0.36%  "
Set-TextCell $ws "D33" "6.399"
Set-TextCell $ws "E33" "  +3.70%  "
Set-TextCell $ws "D34" "1.633"
Set-TextCell $ws "E34" "  -2.26%  "
Set-TextCell $ws "E35" "  +0.22%  "
Set-TextCell $ws "D36" "10.52"
Set-TextCell $ws "E36" "  -0.54%  "
Set-TextCell $ws "D37" "5.943"
Set-TextCell $ws "E37" "  +7.84%  "
Set-TextCell $ws "D38" "0.02671"
Set-TextCell $ws "E38" "  +3.84%  "
Set-TextCell $ws "D39" "0.06861"
Set-TextCell $ws "E39" "  +0.94%  "
Set-TextCell $ws "D40" "0.2323"
Set-TextCell $ws "E40" "  +1.96%  "
Set-TextCell $ws "D41" "12.62"
Set-TextCell $ws "E41" "  -1.17%  "
Set-TextCell $ws "D42" "0.6890"
Set-TextCell $ws "E42" "  -0.67%  "
Set-TextCell $ws "D43" "1.264"
Set-TextCell $ws "E43" "  +0.57%  "
Set-TextCell $ws "D44" "15.22"
Set-TextCell $ws "E44" "  +7.84%  "
Set-TextCell $ws "D45" "0.6463"
Set-TextCell $ws "E45" "  +1.04%  "
Set-TextCell $ws "D46" "2.308"
Set-TextCell $ws "E46" "  -1.88%  "
Set-TextCell $ws "E47" "  +15.46%  "
Set-TextCell $ws "D48" "3.701"
Set-TextCell $ws "E48" "  +1.36%  "
Set-TextCell $ws "D49" "1.256"
Set-TextCell $ws "E49" "  +0.73%  "

# Row 50/51: source re-ranked Cronos above Aave, so the two rows'
# full contents (name/link/price/volume) swap place.
Set-TextCell $ws "B50" "Cronos"
Set-TextCell $ws "C50" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D50" "0.07318"
Set-TextCell $ws "E50" "  +3.51%  "
Set-TextCell $ws "B51" "Aave"
Set-TextCell $ws "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws "D51" "82.89"
Set-TextCell $ws "E51" "  -0.56%  "
